function Set-RowData($ws, $r, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $hVal, $iVal) {
    $ws.Cells.Item($r, 2).Value2 = "'" + $bVal
    $ws.Cells.Item($r, 3).Value2 = $cVal
    $ws.Cells.Item($r, 4).Value2 = $dVal
    $ws.Cells.Item($r, 5).Value2 = $eVal
    $ws.Cells.Item($r, 6).Value2 = $fVal
    $ws.Cells.Item($r, 7).Value2 = $gVal
    $ws.Cells.Item($r, 8).Value2 = $hVal
    $ws.Cells.Item($r, 9).Value2 = $iVal
}

function Update-Sheet($ws) {
    # Remove the last two rows (old rows 37 and 38); rows above are
    # untouched by this operation because it only affects rows below them.
    $ws.Range("A37:I38").EntireRow.Delete()

    Set-RowData $ws 2 '2024-04-20' '赣州·十万伏特-次元音乐only' '平安大道 麋鹿LiveHouse' '2024.04.20 14:30-04.21 21:00' 39 65 'https://show.bilibili.com/platform/detail.html?id=83871' '//i2.hdslb.com/bfs/openplatform/202404/XlyMkr9u1711959548249.jpeg'
    Set-RowData $ws 3 '2024-05-01' '九江·第三届ACD动漫游戏嘉年华' '九瑞大道与重庆路交汇处西南角 九江国际会展中心' '2024.05.01 09:00-05.02 17:00' 689 55 'https://show.bilibili.com/platform/detail.html?id=82464' '//i0.hdslb.com/bfs/openplatform/202403/HjMMyP3a1709780146797.jpeg'
    Set-RowData $ws 4 '2024-05-01' '吉安·COMIC LIFE次元假日04' '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心' '2024.05.01 09:00-05.01 18:00' 686 50 'https://show.bilibili.com/platform/detail.html?id=82949' '//i0.hdslb.com/bfs/openplatform/202403/XDTuzcBV1710478728595.jpeg'
    Set-RowData $ws 5 '2024-05-01' '宜春·静卿璀璨国风动漫文化艺术节' '中山西路 大富豪酒店(国光翰邦店)' '2024.05.01 09:00-05.01 17:00' 39 40 'https://show.bilibili.com/platform/detail.html?id=84514' '//i2.hdslb.com/bfs/openplatform/202404/AHEDWBoB1713111000623.jpeg'
    Set-RowData $ws 6 '2024-05-01' '新余·新次元动漫游戏嘉年华' '渝东大道3208号 新余广富城购物中心' '2024.05.01 10:00-05.04 17:00' 17 30 'https://show.bilibili.com/platform/detail.html?id=84601' '//i1.hdslb.com/bfs/openplatform/202404/yGRxJK4g1713518272823.png'
    Set-RowData $ws 7 '2024-05-01' '景德镇·第一届国际动漫节' '迎宾大道与寺山路交叉口东200米 陶博城' '2024.05.01 10:00-05.02 18:00' 1087 55 'https://show.bilibili.com/platform/detail.html?id=83734' '//i1.hdslb.com/bfs/openplatform/202404/6Izq7ZEk1712054058852.jpeg'
    Set-RowData $ws 8 '2024-05-01' '景德镇·第一届国际动漫节吕书君内场票' '迎宾大道与寺山路交叉口东200米 陶博城' '2024.05.01 10:00-05.01 18:00' 228 128 'https://show.bilibili.com/platform/detail.html?id=83666' '//i1.hdslb.com/bfs/openplatform/202403/OBgwOCB21711786949138.jpeg'
    Set-RowData $ws 9 '2024-05-01' '景德镇·第十四届瓷都ACG动漫游戏博览会' '新厂西路315号 陶溪川发布大厅' '2024.05.01 10:00-05.02 18:00' 1605 50 'https://show.bilibili.com/platform/detail.html?id=83016' '//i2.hdslb.com/bfs/openplatform/202403/c0q8seJL1710835930052.png'
    Set-RowData $ws 10 '2024-05-01' '江西·广电·Unlimited Project 动漫游戏博览会' '怀玉山大道1315号 南昌绿地国际博览中心' '2024.05.01 09:00-05.02 17:00' 5683 68 'https://show.bilibili.com/platform/detail.html?id=83230' '//i2.hdslb.com/bfs/openplatform/202404/S1nqZf721712025221477.jpeg'
    Set-RowData $ws 11 '2024-05-01' '江西·第二十二届九江ACJJ国际动漫展' '体育路九江市体育中心-体育馆 九江市体育中心' '2024.05.01 09:00-05.02 17:00' 466 50 'https://show.bilibili.com/platform/detail.html?id=83004' '//i2.hdslb.com/bfs/openplatform/202403/lFThDvkh1710829330909.jpeg'
    Set-RowData $ws 12 '2024-05-01' '赣州·COMIC WORLD次元创作同人季特典·SP·动漫游戏嘉年华' '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心' '2024.05.01 10:00-05.03 17:00' 310 58 'https://show.bilibili.com/platform/detail.html?id=83537' '//i2.hdslb.com/bfs/openplatform/202403/BKJfMXXx1711091647172.jpeg'
    Set-RowData $ws 13 '2024-05-01' '赣州·十万伏特-第六届青年文化综合展览会' '石楠路与仓背岭路交叉口北120米 新旅中书艺术综合体' '2024.05.01 09:30-05.03 17:00' 255 65 'https://show.bilibili.com/platform/detail.html?id=83628' '//i1.hdslb.com/bfs/openplatform/202404/k8EDbMuk1711961223856.jpeg'
    Set-RowData $ws 14 '2024-05-01' '赣州·超新星国风动漫节' '红旗大道63号 赣坊1969文化创意产业园' '2024.05.01 10:00-05.02 17:00' 67 20 'https://show.bilibili.com/platform/detail.html?id=84419' '//i1.hdslb.com/bfs/openplatform/202404/xSgtoj051713190844041.jpeg'
    Set-RowData $ws 15 '2024-05-02' '南昌·静卿璀璨国风动漫文化艺术节' '京东南大道1666号 万达广场(青山湖店)' '2024.05.02 10:00-05.02 17:00' 37 30 'https://show.bilibili.com/platform/detail.html?id=84516' '//i2.hdslb.com/bfs/openplatform/202404/wVsJVVbN1713114136374.jpeg'
    Set-RowData $ws 16 '2024-05-02' '抚州·临次元07国漫&运动番嘉年华' '伍塘路1098号 乐课篮球公园' '2024.05.02 10:00-05.02 16:00' 121 50 'https://show.bilibili.com/platform/detail.html?id=82911' '//i0.hdslb.com/bfs/openplatform/202403/Pc5TMZ001710429899787.jpeg'
    Set-RowData $ws 17 '2024-05-02' '江西·ShiningStaR数字互娱嘉年华' '前湖大道欣悦湖体育馆 欣悦湖体育馆' '2024.05.02 09:30-05.04 17:00' 4628 65 'https://show.bilibili.com/platform/detail.html?id=83180' '//i0.hdslb.com/bfs/openplatform/202403/EqmGU5NC1711015780862.jpeg'
    Set-RowData $ws 18 '2024-05-03' '上饶·第一届星光次元国风动漫游戏嘉年华暨我和我的cos小伙伴们' '带湖路66-68号 华都臻悦酒店' '2024.05.03 09:30-05.03 17:30' 239 45 'https://show.bilibili.com/platform/detail.html?id=83478' '//i2.hdslb.com/bfs/openplatform/202403/GFCaK00i1711614562687.jpeg'
    Set-RowData $ws 19 '2024-05-03' '南昌·第二届龙年动漫展' '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆' '2024.05.03 10:00-05.04 18:00' 1223 55 'https://show.bilibili.com/platform/detail.html?id=83954' '//i1.hdslb.com/bfs/openplatform/202404/eECmAU3Q1713511863782.jpeg'
    Set-RowData $ws 20 '2024-05-03' '新余·LD02国风动漫嘉年华' '劳动北路888号 金联体育篮球馆' '2024.05.03 10:00-05.03 17:00' 127 50 'https://show.bilibili.com/platform/detail.html?id=83379' '//i0.hdslb.com/bfs/openplatform/202403/jozduadT1711362183223.jpeg'
    Set-RowData $ws 21 '2024-05-03' '江西·2024南昌玛雅《次元之芯》主题动漫嘉年华' '真君路999号 南昌玛雅乐园' '2024.05.03 10:30-05.04 19:30' 87 56 'https://show.bilibili.com/platform/detail.html?id=83608' '//i1.hdslb.com/bfs/openplatform/202404/n9Ba2rOM1713338456905.jpeg'
    Set-RowData $ws 22 '2024-05-03' '江西·ShiningStaR数字互娱嘉年华配音演员史泽鲲专场见面会' '前湖大道欣悦湖体育馆 欣悦湖体育馆' '2024.05.03 09:30-05.03 17:30' 214 188 'https://show.bilibili.com/platform/detail.html?id=83497' '//i1.hdslb.com/bfs/openplatform/202403/qm19B8RF1711620646864.jpeg'
    Set-RowData $ws 23 '2024-05-03' '萍乡·AU8春季国漫展' '金陵东路18号 萍乡市体育馆' '2024.05.03 10:00-05.03 17:00' 82 45 'https://show.bilibili.com/platform/detail.html?id=84259' '//i0.hdslb.com/bfs/openplatform/202404/GxgrLKhd1712560102907.png'
    Set-RowData $ws 24 '2024-05-03' '赣州·漫库书店次元漫展' '南门口地一大道下沉广场 漫库书店' '2024.05.03 10:00-05.04 18:00' 223 45 'https://show.bilibili.com/platform/detail.html?id=83855' '//i0.hdslb.com/bfs/openplatform/202404/juDVRy6Y1712481590113.jpeg'
    Set-RowData $ws 25 '2024-05-04' '南昌·Youth动漫美食嘉年华' '灌婴路西口朝阳江滩公园内 元亨体育运动中心' '2024.05.04 09:00-05.05 20:00' 77 55 'https://show.bilibili.com/platform/detail.html?id=84036' '//i2.hdslb.com/bfs/openplatform/202404/I5vd5js01712648032400.jpeg'
    Set-RowData $ws 26 '2024-05-04' '江西·ShiningStaR数字互娱嘉年华 配音演员陈张太康、张惠霖专场见面会' '前湖大道欣悦湖体育馆 欣悦湖体育馆' '2024.05.04 09:30-05.04 17:30' 156 228 'https://show.bilibili.com/platform/detail.html?id=83593' '//i0.hdslb.com/bfs/openplatform/202404/LcnCzDxF1711935576170.jpeg'
    Set-RowData $ws 27 '2024-05-18' '南昌·花绒万兽首届兽聚' '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆' '2024.05.18 09:30-05.19 16:30' 85 60 'https://show.bilibili.com/platform/detail.html?id=83689' '//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg'
    Set-RowData $ws 28 '2024-05-26' '南昌·代号鸢盛花行only' '民德路411号 东方豪景花园酒店(民德路店)' '2024.05.26 09:30-05.26 17:30' 360 78 'https://show.bilibili.com/platform/detail.html?id=82529' '//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png'
    Set-RowData $ws 29 '2024-06-01' '南昌·ACG CLUB动漫游戏嘉年华' '南京东路1225号橙风体育(青山湖店) 橙风体育篮球馆(青山湖店)' '2024.06.01 10:00-06.02 17:00' 45 55 'https://show.bilibili.com/platform/detail.html?id=84497' '//i2.hdslb.com/bfs/openplatform/202404/PGGkrHaZ1712914809144.jpeg'
    Set-RowData $ws 30 '2024-06-09' '信丰·端午节UPUP动漫展' '迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子体育馆' '2024.06.09 10:00-06.09 17:00' 48 48 'https://show.bilibili.com/platform/detail.html?id=84078' '//i0.hdslb.com/bfs/openplatform/202404/Qy0EOl551712651477492.jpeg'
    Set-RowData $ws 31 '2024-06-10' '上饶·ETI动漫节' '滨江东路与体育馆路交叉口西100米 力加体育综合运动中心' '2024.06.10 10:00-06.10 16:00' 71 36.6 'https://show.bilibili.com/platform/detail.html?id=83422' '//i2.hdslb.com/bfs/openplatform/202404/N6VdMOuL1713257425864.jpeg'
    Set-RowData $ws 32 '2024-06-10' '南昌·LY-COSPLAY大会X运动番PRO2.0（非ONLY）' '青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK' '2024.06.10 10:00-06.10 17:00' 24 30 'https://show.bilibili.com/platform/detail.html?id=84575' '//i2.hdslb.com/bfs/openplatform/202404/ScwkijwU1713428452963.jpeg'
    Set-RowData $ws 33 '2024-08-03' '南昌·幻梦境国际动漫游戏嘉年华1th' '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆' '2024.08.03 09:00-08.04 17:30' 36 64 'https://show.bilibili.com/platform/detail.html?id=83980' '//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'
    Set-RowData $ws 34 '2024-08-04' '九江·第一届异次元动漫嘉年华' '长虹西大道兴城广场99号 九江半岛宾馆' '2024.08.04 08:00-08.04 17:00' 9 35 'https://show.bilibili.com/platform/detail.html?id=84407' '//i1.hdslb.com/bfs/openplatform/202404/e7k26XLV1713262153782.jpeg'
    Set-RowData $ws 35 '2024-08-06' '南昌·第一届异次元动漫嘉年华' '民德路411号 东方豪景花园酒店(民德路店)' '2024.08.06 08:00-08.06 17:00' 49 40 'https://show.bilibili.com/platform/detail.html?id=84102' '//i0.hdslb.com/bfs/openplatform/202404/0W8gIOYx1712764727806.jpeg'
    Set-RowData $ws 36 '2024-08-08' '赣州·第二届异次元动漫嘉年华' '金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆' '2024.08.08 08:00-08.08 17:00' 54 45 'https://show.bilibili.com/platform/detail.html?id=84184' '//i2.hdslb.com/bfs/openplatform/202404/F5F9vvqX1712758945373.jpeg'

}

$wb = $excel.ActiveWorkbook

Update-Sheet $wb.Worksheets.Item(1)
Update-Sheet $wb.Worksheets.Item(4)

Write-Output "Done"
